$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 4.194610486936591
$ws.Range("D2").Value = 9.569347722294449
$ws.Range("E2").Value = 13.84292599521781
$ws.Range("F2").Value = 31.23611681445975
$ws.Range("G2").Value = 3.64737526039056
$ws.Range("I2").Value = 21.6320821660381
$ws.Range("J2").Value = 9.872518522162169
$ws.Range("K2").Value = 17.3022849943482
$ws.Range("O2").Value = 23.39063903491109
$ws.Range("C3").Value = 4.080754817109423
$ws.Range("D3").Value = 9.486150388326237
$ws.Range("E3").Value = 13.76896597225797
$ws.Range("F3").Value = 31.37544491137397
$ws.Range("G3").Value = 3.649657766113261
$ws.Range("I3").Value = 21.7536717553161
$ws.Range("J3").Value = 9.881908553824912
$ws.Range("K3").Value = 16.60921730026283
$ws.Range("O3").Value = 23.54478484021562
$ws.Range("C4").Value = 4.008639872439859
$ws.Range("D4").Value = 9.435861461953655
$ws.Range("E4").Value = 13.72604734753341
$ws.Range("F4").Value = 31.47134251587159
$ws.Range("G4").Value = 3.651131723957966
$ws.Range("I4").Value = 21.83513027385881
$ws.Range("J4").Value = 9.889341373462049
$ws.Range("K4").Value = 16.16834046132578
$ws.Range("O4").Value = 23.6466012201505
$ws.Range("C5").Value = 3.978725405040316
$ws.Range("D5").Value = 9.415585151080871
$ws.Range("E5").Value = 13.70919814481828
$ws.Range("F5").Value = 31.51301000273679
$ws.Range("G5").Value = 3.6517506613563
$ws.Range("I5").Value = 21.87002796843324
$ws.Range("J5").Value = 9.892789555397641
$ws.Range("K5").Value = 15.98504849944392
$ws.Range("O5").Value = 23.68988877388743
$ws.Range("C6").Value = 3.97372715782322
$ws.Range("D6").Value = 9.412231871653997
$ws.Range("E6").Value = 13.70643940995901
$ws.Range("F6").Value = 31.52008477404822
$ws.Range("G6").Value = 3.651854541712664
$ws.Range("I6").Value = 21.87592531577196
$ws.Range("J6").Value = 9.893387442600208
$ws.Range("K6").Value = 15.95440079643872
$ws.Range("O6").Value = 23.69718496002385
$ws.Range("C7").Value = 4.008238530926345
$ws.Range("D7").Value = 9.4355871083502
$ws.Range("E7").Value = 13.72581750273441
$ws.Range("F7").Value = 31.47189399661017
$ws.Range("G7").Value = 3.651139997039631
$ws.Range("I7").Value = 21.83559403253041
$ws.Range("J7").Value = 9.889386179423175
$ws.Range("K7").Value = 16.1658829074648
$ws.Range("O7").Value = 23.64717774556613
$ws.Range("C8").Value = 4.15582723303988
$ws.Range("D8").Value = 9.540506694133462
$ws.Range("E8").Value = 13.81691553894355
$ws.Range("F8").Value = 31.28200075062732
$ws.Range("G8").Value = 3.648147258878002
$ws.Range("I8").Value = 21.67258985266911
$ws.Range("J8").Value = 9.875410179419077
$ws.Range("K8").Value = 17.06662289652471
$ws.Range("O8").Value = 23.44229641289073
$ws.Range("C9").Value = 4.42661872121478
$ws.Range("D9").Value = 9.751795429294111
$ws.Range("E9").Value = 14.01470435654661
$ws.Range("F9").Value = 30.99235056580243
$ws.Range("G9").Value = 3.642850965887483
$ws.Range("I9").Value = 21.40725493987252
$ws.Range("J9").Value = 9.861229922837545
$ws.Range("K9").Value = 18.70285628469509
$ws.Range("O9").Value = 23.0977101675796
$ws.Range("C10").Value = 4.612868529172054
$ws.Range("D10").Value = 9.90937096609378
$ws.Range("E10").Value = 14.17078646728502
$ws.Range("F10").Value = 30.830785855673
$ws.Range("G10").Value = 3.639304945569891
$ws.Range("I10").Value = 21.24589570712284
$ws.Range("J10").Value = 9.858866795763797
$ws.Range("K10").Value = 19.81605735060812
$ws.Range("O10").Value = 22.87981161379226
$ws.Range("C11").Value = 4.694597892708351
$ws.Range("D11").Value = 9.981341663754224
$ws.Range("E11").Value = 14.243922604925
$ws.Range("F11").Value = 30.76857502406443
$ws.Range("G11").Value = 3.637765904245766
$ws.Range("I11").Value = 21.17988151057968
$ws.Range("J11").Value = 9.859536842505033
$ws.Range("K11").Value = 20.30157559152633
$ws.Range("O11").Value = 22.78843261815602
$ws.Range("C12").Value = 4.725098083487429
$ws.Range("D12").Value = 10.00861786132424
$ws.Range("E12").Value = 14.2719057577
$ws.Range("F12").Value = 30.74665278436012
$ws.Range("G12").Value = 3.637193697523984
$ws.Range("I12").Value = 21.15595443037151
$ws.Range("J12").Value = 9.860040936978203
$ws.Range("K12").Value = 20.4823178996291
$ws.Range("O12").Value = 22.75495120886878
$ws.Range("C13").Value = 4.718549556343587
$ws.Range("D13").Value = 10.00274274765844
$ws.Range("E13").Value = 14.26586657563768
$ws.Range("F13").Value = 30.75130118830296
$ws.Range("G13").Value = 3.637316462144645
$ws.Range("I13").Value = 21.16105978405225
$ws.Range("J13").Value = 9.859921246730641
$ws.Range("K13").Value = 20.44353183612653
$ws.Range("O13").Value = 22.76211200290178
$ws.Range("C14").Value = 4.697116256076953
$ws.Range("D14").Value = 9.983585351664841
$ws.Range("E14").Value = 14.24621913105424
$ws.Range("F14").Value = 30.76673861642743
$ws.Range("G14").Value = 3.63771861643246
$ws.Range("I14").Value = 21.1778915008916
$ws.Range("J14").Value = 9.859573300751268
$ws.Range("K14").Value = 20.31650832953216
$ws.Range("O14").Value = 22.78565555157293
$ws.Range("C15").Value = 4.683928770969318
$ws.Range("D15").Value = 9.971853254069465
$ws.Range("E15").Value = 14.23422145710169
$ws.Range("F15").Value = 30.77640786524497
$ws.Range("G15").Value = 3.63796632567629
$ws.Range("I15").Value = 21.18834114994201
$ws.Range("J15").Value = 9.859392759227834
$ws.Range("K15").Value = 20.23829438566114
$ws.Range("O15").Value = 22.80022299761005
$ws.Range("C16").Value = 4.607465439447254
$ws.Range("D16").Value = 9.904671690085278
$ws.Range("E16").Value = 14.16604819292084
$ws.Range("F16").Value = 30.83507968457223
$ws.Range("G16").Value = 3.639407011122221
$ws.Range("I16").Value = 21.25035930614984
$ws.Range("J16").Value = 9.858858075251245
$ws.Range("K16").Value = 19.7838975957678
$ws.Range("O16").Value = 22.88593997863277
$ws.Range("C17").Value = 4.559776733347336
$ws.Range("D17").Value = 9.863518340798368
$ws.Range("E17").Value = 14.12475940857292
$ws.Range("F17").Value = 30.87397278774442
$ws.Range("G17").Value = 3.640309754506368
$ws.Range("I17").Value = 21.29030401782913
$ws.Range("J17").Value = 9.858976643512349
$ws.Range("K17").Value = 19.49970772619966
$ws.Range("O17").Value = 22.9405138442473
$ws.Range("C18").Value = 4.53206635263635
$ws.Range("D18").Value = 9.839876342163075
$ws.Range("E18").Value = 14.10121317168722
$ws.Range("F18").Value = 30.89740461614715
$ws.Range("G18").Value = 3.640835962873944
$ws.Range("I18").Value = 21.31397413290339
$ws.Range("J18").Value = 9.859209134865676
$ws.Range("K18").Value = 19.33428911586197
$ws.Range("O18").Value = 22.97263170237564
$ws.Range("C19").Value = 4.522636395020802
$ws.Range("D19").Value = 9.831877014448606
$ws.Range("E19").Value = 14.09327607298839
$ws.Range("F19").Value = 30.90552015243828
$ws.Range("G19").Value = 3.641015327504278
$ws.Range("I19").Value = 21.3221075281525
$ws.Range("J19").Value = 9.859316085531875
$ws.Range("K19").Value = 19.27794819268038
$ws.Range("O19").Value = 22.98363111469071
$ws.Range("C20").Value = 4.564882506657036
$ws.Range("D20").Value = 9.867896392669921
$ws.Range("E20").Value = 14.12913390008143
$ws.Range("F20").Value = 30.8697225900921
$ws.Range("G20").Value = 3.6402129345232
$ws.Range("I20").Value = 21.28597984579488
$ws.Range("J20").Value = 9.85894702165945
$ws.Range("K20").Value = 19.53016393695875
$ws.Range("O20").Value = 22.93462891574103
$ws.Range("C21").Value = 4.703424054799165
$ws.Range("D21").Value = 9.989211885141422
$ws.Range("E21").Value = 14.25198239485602
$ws.Range("F21").Value = 30.76215977983393
$ws.Range("G21").Value = 3.637600206866795
$ws.Range("I21").Value = 21.17291848154202
$ws.Range("J21").Value = 9.85966871134155
$ws.Range("K21").Value = 20.35390349189684
$ws.Range("O21").Value = 22.77870972804182
$ws.Range("C22").Value = 4.791345075144528
$ws.Range("D22").Value = 10.06862072687564
$ws.Range("E22").Value = 14.33394116108645
$ws.Range("F22").Value = 30.70140056088838
$ws.Range("G22").Value = 3.635954367760212
$ws.Range("I22").Value = 21.1052732332961
$ws.Range("J22").Value = 9.861599326671872
$ws.Range("K22").Value = 20.87408058687837
$ws.Range("O22").Value = 22.68335035159236
$ws.Range("C23").Value = 4.744665546962926
$ws.Range("D23").Value = 10.02623381983263
$ws.Range("E23").Value = 14.29005170967
$ws.Range("F23").Value = 30.73295200139086
$ws.Range("G23").Value = 3.636827152919104
$ws.Range("I23").Value = 21.1408023802258
$ws.Range("J23").Value = 9.860435648440912
$ws.Range("K23").Value = 20.59814835761903
$ws.Range("O23").Value = 22.73364401821913
$ws.Range("C24").Value = 4.562575099380662
$ws.Range("D24").Value = 9.865917019131922
$ws.Range("E24").Value = 14.12715559539621
$ws.Range("F24").Value = 30.87164076718372
$ws.Range("G24").Value = 3.640256684376928
$ws.Range("I24").Value = 21.28793260738247
$ws.Range("J24").Value = 9.858959901799052
$ws.Range("K24").Value = 19.51640101124758
$ws.Range("O24").Value = 22.93728717925812
$ws.Range("C25").Value = 4.355503752064717
$ws.Range("D25").Value = 9.694149326682092
$ws.Range("E25").Value = 13.95923801074832
$ws.Range("F25").Value = 31.06176658858483
$ws.Range("G25").Value = 3.644222864098686
$ws.Range("I25").Value = 21.47317155601469
$ws.Range("J25").Value = 9.863650444611228
$ws.Range("K25").Value = 18.27522318448346
$ws.Range("O25").Value = 23.1847663000685
